# Fill in the "FORMTEXT" field that follows the "...circumstances and
# behaviours leading to the recall..." heading with the merge-field
# placeholder {{what_led_to_recall}}.
#
# The field's result is currently five runs each holding a single
# space character; we replace the whole result range in one shot so it
# collapses to a single run containing the placeholder text.

$d = $word.ActiveDocument

$headingIndex = -1
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text
    if ($paraText -match "circumstances and behaviours leading to the recall") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    Write-Output "Heading paragraph not found"
} else {
    $fieldPara = $d.Paragraphs.Item($headingIndex + 1)
    $fieldParaRange = $fieldPara.Range

    $replaced = $false
    foreach ($fld in $d.Fields) {
        if ($fld.Type -eq 1 -and $fld.Code.Text -match "FORMTEXT") {
            $resultRange = $fld.Result
            if ($resultRange.Start -ge $fieldParaRange.Start -and $resultRange.End -le $fieldParaRange.End) {
                $targetRange = $d.Range($resultRange.Start, $resultRange.End)
                $targetRange.Text = "{{what_led_to_recall}}"
                $replaced = $true
                Write-Output "Updated FORMTEXT field result to: $($targetRange.Text)"
                break
            }
        }
    }

    if (-not $replaced) {
        Write-Output "Target FORMTEXT field not found"
    }
}
